$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Switch to single faculty list: row 11 (Heid, Leslie) now also has an
# asterisk flag set, and its faculty-name value moves from the combined
# "Allard, Lander" entry to the single name "Lander".
$ws.Range("C11").Value = 1
$ws.Range("D11").Value = "Lander"

# Move the active selection to D11 (was D20).
$ws.Range("D11").Select()
